$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "book" references to "story": BCH_xxx -> SCH_xxx in the Links cell (E2)
$ws.Range("E2").Value = "CH_004, CH_005, CH_006, CH_007, CH_008, CH_011, CH_018,SCH_003, SCH_005, SCH_006, SCH_007, SCH_008, SCH_009, SCH_010, SCH_011, SCH_012, SCH_013, SCH_019, SCH_020, SCH_021, SCH_022, SCH_023, SCH_024, SCH_025"

# Select the entire column E (active cell E1), matching the saved selection state
$ws.Columns("E").Select()
